$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.InsertBefore("Update") | Out-Null
$s.Shapes.Item(2).TextFrame.TextRange.LanguageID = "en-US"
